$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Date
$ws.Range("B3").Value = "12/2/2020"

# Row 4: Team Name
$ws.Range("B4").Value = "Limette"

# Row 5: Total Number of Team Members
$ws.Range("B5").Value = 4

# Rows 8-11: team member names and salaries
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 100

$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 100

$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 100

$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 100

# Row 12: clear leftover "Member 5" label
$ws.Range("A12").Value = ""

# Row 18: shrink the row so the two header cells fit tighter
$ws.Range("A18").RowHeight = 39

# Rows 19-20: new task notes using a smaller 10pt Calibri font (matches style added to styles.xml)
$ws.Range("A19").Font.Size = 10
$ws.Range("B19").Font.Size = 10
$ws.Range("A20").Font.Size = 10
$ws.Range("B20").Font.Size = 10

$ws.Range("A19").Value = "Brainstorming for optimization possibilities."
$ws.Range("B19").Value = "Finish automated hifi prototype."
$ws.Range("A20").Value = "Started working on the automated hifi prototype."

# Update the selected cell to match the saved session state
[void]$ws.Range("E16").Select()
